$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Pages")

# Row 33 updates
$ws.Range("A33").Value = 326
$ws.Range("B33").Value = "Tale"
$ws.Range("C33").Value = 744
$ws.Range("D33").Value = 4631
$ws.Range("E33").Value = 0.9529259339235586
$ws.Range("F33").Value = 5.931451612903226
$ws.Range("G33").Value = 45801.52984953704

# Row 41 updates
$ws.Range("C41").Value = 165216
$ws.Range("D41").Value = 1058861
$ws.Range("E41").Value = 0.976415223527923
$ws.Range("F41").Value = 6.257795855122991
$ws.Range("G41").Value = 45801.56623842593

# Row 58 updates
$ws.Range("B58").Value = "Tale"
$ws.Range("D58").Value = 13886
$ws.Range("E58").Value = 0.988693648278842
$ws.Range("G58").Value = 45801.53
